$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.675.92'
$ws.Range("E2").Value = '  -2.46%  '

$ws.Range("D3").Value = '3.275.75'
$ws.Range("E3").Value = '  -4.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -10.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").Value = '3.271.09'
$ws.Range("E8").Value = '  -4.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -11.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.61'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.505'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -13.02%  '

$ws.Range("E14").Value = '  -8.59%  '

$ws.Range("D15").Value = '3.806.76'
$ws.Range("E15").Value = '  -4.76%  '

$ws.Range("D16").Value = '67.751.68'
$ws.Range("E16").Value = '  -2.33%  '

$ws.Range("D17").Value = '3.281.25'
$ws.Range("E17").Value = '  -4.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '531.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.59%  '

$ws.Range("E19").Value = '  -5.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -12.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -12.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.756'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -11.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -10.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.91%  '

$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.19%  '

$ws.Range("E31").Value = '  -2.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.39%  '

$ws.Range("E33").Value = '  -15.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -12.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '56.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.69%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '514.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0448'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.63%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0857'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -14.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.126'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.46%  '

$ws.Range("D43").Value = '2.940.04'
$ws.Range("E43").Value = '  -9.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.268'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.09%  '

$ws.Range("D45").Value = '0.0₃0587'
$ws.Range("E45").Value = '  -14.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -14.17%  '

$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -15.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.113'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.27%  '

